$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.492611169815063
$ws.Range("B1").Value = 6.425615787506104
$ws.Range("C1").Value = 3.518653392791748
$ws.Range("D1").Value = 1.571729183197021
$ws.Range("E1").Value = 1.107196927070618
